$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.210.41"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "2.494.81"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +3.19%  "
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "2.885.28"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "2.500.53"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.852"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "47.141.24"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.141"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0786"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("E38").Value = "  +2.53%  "
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0298"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").Value = "1.993.91"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.48%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.13%  "
